$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "30.683.26"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  -0.62%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.889.06"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  -0.78%  "

$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  +0.19%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "236.78"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  -3.99%  "

$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  +0.24%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4879"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  -2.60%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2917"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  -2.79%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06679"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  -2.67%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.890.79"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  -0.62%  "

$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  -3.92%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.07256"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  -1.35%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "89.24"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  -2.76%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.022"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  -2.12%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.6618"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  -3.06%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "30.621.75"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  -0.79%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.000007881"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  -2.22%  "

$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  +0.24%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "12.99"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  -2.82%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "2.135.26"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  -0.79%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "1.003"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  +0.08%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.742"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  -2.91%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "190.50"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  +3.28%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "6.112"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  -0.11%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "9.315"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  -0.86%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "159.12"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  +3.13%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "18.30"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  -1.73%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.835"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  -6.17%  "

$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  +0.59%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.252"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  -3.32%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.09008"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  +0.14%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.933"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  -3.66%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.05170"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  -2.24%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.7258"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  -2.81%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.083"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  -5.46%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.694"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  +0.94%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.01809"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  -6.24%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.665"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  -2.22%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.9218"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  -2.27%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.041"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  -7.03%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.4386"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  -0.45%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "104.53"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  -1.71%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.9999"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  -0.06%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "5.724"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  -2.41%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.1328"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  -2.20%  "

$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  -6.09%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.4058"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  +3.04%  "

$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  -0.40%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "8.668"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  +0.83%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.413"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  +0.99%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "33.30"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  -0.58%  "
